{"js": "// Applies the LOQ4261 course-sheet update:\n//  - Course title (EN) gets a trailing \" I\"\n//  - Ativa\u00e7\u00e3o date 2021 -> 2024\n//  - \"Programa resumido\"/\"Programa\" summaries: item 2 (demand mgmt / gest\u00e3o e\n//    previs\u00e3o de demanda) removed and a new item 8 (Tambor-Pulm\u00e3o-Corda/\n//    Drum-Buffer-Rope) is added; list renumbered accordingly\n//  - M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o rewritten in \"Avalia\u00e7\u00e3o\"\n\nconst body = context.document.body;\n\n// Replace the occurrenceIndex-th (0-based, in document order) match of\n// searchText with newText. Re-searches every call since earlier replacements\n// shift later occurrence indices down by one.\nasync function replaceOccurrence(searchText, newText, occurrenceIndex) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length <= occurrenceIndex) {\n    throw new Error(\n      \"Expected occurrence \" + occurrenceIndex + \" of '\" +\n      searchText.substring(0, 60) + \"...' but only found \" + results.items.length\n    );\n  }\n\n  results.items[occurrenceIndex].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nasync function replaceFirst(searchText, newText) {\n  await replaceOccurrence(searchText, newText, 0);\n}\n\n// 1) Heading3 course title (English) \u2014 unique in the document.\nawait replaceFirst(\n  \"Production Planning, Scheduling and Control\",\n  \"Production Planning, Scheduling and Control I\"\n);\n\n// 2) Ativa\u00e7\u00e3o date \u2014 unique in the document.\nawait replaceFirst(\n  \"Ativa\u00e7\u00e3o: 01/01/2021\",\n  \"Ativa\u00e7\u00e3o: 01/01/2024\"\n);\n\n// 3) Portuguese programme summary \u2014 identical text appears twice\n//    (\"Programa resumido\" then \"Programa\"); both get the same new text.\nconst ptOld =\n  \"1. Caracteriza\u00e7\u00e3o do planejamento e controle da produ\u00e7\u00e3o. 2. Gest\u00e3o e \" +\n  \"previs\u00e3o de demanda. 3. Planejamento agregado da produ\u00e7\u00e3o. 4. \" +\n  \"Planejamento mestre da produ\u00e7\u00e3o. 5. Planejamento e controle de \" +\n  \"estoques. 6. Planejamento de recursos de materiais (MRP). 7. \" +\n  \"Programa\u00e7\u00e3o detalhada da produ\u00e7\u00e3o. 8. Just In Time (JIT) e opera\u00e7\u00f5es \" +\n  \"enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restri\u00e7\u00f5es \" +\n  \"(TOC). 11. Sistemas de controle da produ\u00e7\u00e3o.\";\nconst ptNew =\n  \"Caracteriza\u00e7\u00e3o do planejamento e controle da produ\u00e7\u00e3o. 2. Planejamento \" +\n  \"agregado da produ\u00e7\u00e3o. 3. Planejamento mestre da produ\u00e7\u00e3o. 4. \" +\n  \"Planejamento e controle de estoques. 5. Planejamento de recursos de \" +\n  \"materiais (MRP). 6. Programa\u00e7\u00e3o detalhada da produ\u00e7\u00e3o. 7. Sistema \" +\n  \"MRPII e Sistema ERP. 8.Tambor-Pulm\u00e3o-Corda - OPT. 9. Teoria das \" +\n  \"Restri\u00e7\u00f5es (TOC).\";\n\nawait replaceOccurrence(ptOld, ptNew, 0); // \"Programa resumido\" occurrence\nawait replaceOccurrence(ptOld, ptNew, 0); // \"Programa\" occurrence (now index 0)\n\n// 4) English programme summary (italic) \u2014 same source text appears twice,\n//    but the two replacements differ: the first occurrence (under\n//    \"Programa resumido\") drops the leading \"1. \", the second occurrence\n//    (under \"Programa\") keeps it.\nconst enOld =\n  \"1. Characterization of production programming and control. 2. Demand \" +\n  \"management. 3. Aggregate Production Planning. 4. Master Production \" +\n  \"Schedulling. 5. Inventory planning and control. 6. Material \" +\n  \"Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed \" +\n  \"scheduling of production. 9. Just In Time (JIT). 10. Theory of \" +\n  \"Constraints (TOC). 11. Production control systems.\";\nconst enNewNoLeadingNumber =\n  \"Characterization of production programming and control. 2. Aggregate \" +\n  \"Production Planning. 3. Master Production Schedulling. 4. Inventory \" +\n  \"planning and control. 5. Material Requirement Planning (MRP). 6. \" +\n  \"Detailed scheduling of production. 7. Production control systems. 8. \" +\n  \"Drum-Buffer-Rope \\u2013 Opt; 9. Theory of Constraints (TOC)\";\nconst enNewWithLeadingNumber =\n  \"1. Characterization of production programming and control. 2. \" +\n  \"Aggregate Production Planning. 3. Master Production Schedulling. 4. \" +\n  \"Inventory planning and control. 5. Material Requirement Planning \" +\n  \"(MRP). 6. Detailed scheduling of production. 7. Production control \" +\n  \"systems. 8. Drum-Buffer-Rope \\u2013 Opt; 9. Theory of Constraints (TOC)\";\n\nawait replaceOccurrence(enOld, enNewNoLeadingNumber, 0);   // \"Programa resumido\" occurrence\nawait replaceOccurrence(enOld, enNewWithLeadingNumber, 0); // \"Programa\" occurrence (now index 0)\n\n// 5) Avalia\u00e7\u00e3o \u2014 M\u00e9todo\nawait replaceFirst(\n  \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios. MANTIDO\",\n  \"Provas, atividades em grupo e atividades individuais.\"\n);\n\n// 6) Avalia\u00e7\u00e3o \u2014 Crit\u00e9rio\nawait replaceFirst(\n  \"M = (0,8P + 0,2T)P = m\u00e9dia aritm\u00e9tica de duas provas escritasT = \" +\n  \"M\u00e9dia das notas de trabalhos e exerc\u00edciosM = M\u00e9dia de aproveitamento \" +\n  \"do alunoAprova\u00e7\u00e3o com m\u00e9dia de aproveitamento maior ou igual a 5,0 e \" +\n  \"no m\u00ednimo 70% de frequ\u00eancia \u00e0s aulas.\",\n  \"M\u00e9dia das atividades avaliativas\"\n);\n\n// 7) Avalia\u00e7\u00e3o \u2014 Norma de recupera\u00e7\u00e3o\nawait replaceFirst(\n  \"MF = (0,5 M + 0,5 R)M = M\u00e9dia de aproveitamento do aluno, antes da \" +\n  \"recupera\u00e7\u00e3oR = Nota de uma prova de recupera\u00e7\u00e3oMF = nota final de \" +\n  \"aproveitamento, ap\u00f3s a recupera\u00e7\u00e3oAprova\u00e7\u00e3o com m\u00e9dia final de \" +\n  \"aproveitamento maior ou igual a 5,0.A recupera\u00e7\u00e3o dever\u00e1 consistir \" +\n  \"de uma prova escrita englobando a mat\u00e9ria toda do semestre.Ter\u00e1 \" +\n  \"direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota \" +\n  \"acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\",\n  \"MF = (0,5 M + 0,5 R) M = M\u00e9dia de aproveitamento do aluno, antes da \" +\n  \"recupera\u00e7\u00e3o R = Nota de uma prova de recupera\u00e7\u00e3o MF = nota final de \" +\n  \"aproveitamento, ap\u00f3s a recupera\u00e7\u00e3o Aprova\u00e7\u00e3o com m\u00e9dia final de \" +\n  \"aproveitamento maior ou igual a 5,0. A recupera\u00e7\u00e3o dever\u00e1 consistir \" +\n  \"em uma prova escrita englobando a mat\u00e9ria toda do semestre. Ter\u00e1 \" +\n  \"direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota \" +\n  \"acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\"\n);\n", "ps1": "# Applies the LOQ4261 course-sheet update:\n#  - Course title (EN) gets a trailing \" I\"\n#  - Ativa\u00e7\u00e3o date 2021 -> 2024\n#  - \"Programa resumido\"/\"Programa\" summaries: item 2 (demand mgmt / gest\u00e3o e\n#    previs\u00e3o de demanda) removed and a new item 8 (Tambor-Pulm\u00e3o-Corda /\n#    Drum-Buffer-Rope) is added; list renumbered accordingly\n#  - M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o rewritten in \"Avalia\u00e7\u00e3o\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange($Range, $FindText, $ReplaceText) {\n    $Range.Find.Execute(\n        $FindText,\n        $false,  # MatchCase\n        $true,   # MatchWholeWord\n        $false,  # MatchWildcards\n        $false,  # MatchSoundsLike\n        $false,  # MatchAllWordForms\n        $true,   # Forward\n        1,       # Wrap (wdFindContinue)\n        $false,  # Format\n        $ReplaceText,\n        2        # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n\n# 1) Heading3 course title (English) \u2014 unique in the document.\nReplace-InRange $d.Paragraphs.Item(2).Range `\n    \"Production Planning, Scheduling and Control\" `\n    \"Production Planning, Scheduling and Control I\"\n\n# 2) Ativa\u00e7\u00e3o date \u2014 inside the bullet list paragraph with Cr\u00e9ditos/Ativa\u00e7\u00e3o/etc.\nReplace-InRange $d.Paragraphs.Item(4).Range `\n    \"Ativa\u00e7\u00e3o: 01/01/2021\" `\n    \"Ativa\u00e7\u00e3o: 01/01/2024\"\n\n# 3) Portuguese programme summary \u2014 identical text appears twice\n#    (\"Programa resumido\" paragraph 11, \"Programa\" paragraph 14); both get\n#    the same new text.\n$ptOld = \"1. Caracteriza\u00e7\u00e3o do planejamento e controle da produ\u00e7\u00e3o. 2. Gest\u00e3o e previs\u00e3o de demanda. 3. Planejamento agregado da produ\u00e7\u00e3o. 4. Planejamento mestre da produ\u00e7\u00e3o. 5. Planejamento e controle de estoques. 6. Planejamento de recursos de materiais (MRP). 7. Programa\u00e7\u00e3o detalhada da produ\u00e7\u00e3o. 8. Just In Time (JIT) e opera\u00e7\u00f5es enxutas. 9. Sistema MRPII e Sistema ERP; 10. Teoria das Restri\u00e7\u00f5es (TOC). 11. Sistemas de controle da produ\u00e7\u00e3o.\"\n$ptNew = \"Caracteriza\u00e7\u00e3o do planejamento e controle da produ\u00e7\u00e3o. 2. Planejamento agregado da produ\u00e7\u00e3o. 3. Planejamento mestre da produ\u00e7\u00e3o. 4. Planejamento e controle de estoques. 5. Planejamento de recursos de materiais (MRP). 6. Programa\u00e7\u00e3o detalhada da produ\u00e7\u00e3o. 7. Sistema MRPII e Sistema ERP. 8.Tambor-Pulm\u00e3o-Corda - OPT. 9. Teoria das Restri\u00e7\u00f5es (TOC).\"\n\nReplace-InRange $d.Paragraphs.Item(11).Range $ptOld $ptNew\nReplace-InRange $d.Paragraphs.Item(14).Range $ptOld $ptNew\n\n# 4) English programme summary (italic) \u2014 same source text appears twice\n#    (paragraph 12 under \"Programa resumido\", paragraph 15 under \"Programa\"),\n#    but the two replacements differ: paragraph 12 drops the leading \"1. \",\n#    paragraph 15 keeps it.\n$enOld = \"1. Characterization of production programming and control. 2. Demand management. 3. Aggregate Production Planning. 4. Master Production Schedulling. 5. Inventory planning and control. 6. Material Requirement Planning (MRP). 7. Production Schedulling. 8. Detailed scheduling of production. 9. Just In Time (JIT). 10. Theory of Constraints (TOC). 11. Production control systems.\"\n$enNewNoLeadingNumber = \"Characterization of production programming and control. 2. Aggregate Production Planning. 3. Master Production Schedulling. 4. Inventory planning and control. 5. Material Requirement Planning (MRP). 6. Detailed scheduling of production. 7. Production control systems. 8. Drum-Buffer-Rope \u2013 Opt; 9. Theory of Constraints (TOC)\"\n$enNewWithLeadingNumber = \"1. Characterization of production programming and control. 2. Aggregate Production Planning. 3. Master Production Schedulling. 4. Inventory planning and control. 5. Material Requirement Planning (MRP). 6. Detailed scheduling of production. 7. Production control systems. 8. Drum-Buffer-Rope \u2013 Opt; 9. Theory of Constraints (TOC)\"\n\nReplace-InRange $d.Paragraphs.Item(12).Range $enOld $enNewNoLeadingNumber\nReplace-InRange $d.Paragraphs.Item(15).Range $enOld $enNewWithLeadingNumber\n\n# 5/6/7) Avalia\u00e7\u00e3o paragraph (17) \u2014 M\u00e9todo / Crit\u00e9rio / Norma de recupera\u00e7\u00e3o\n#    are three separate runs inside the same paragraph; scope each\n#    Find/Replace there so only the intended run's text is touched.\n$avaliacao = $d.Paragraphs.Item(17).Range\n\nReplace-InRange $avaliacao `\n    \"Aulas expositivas te\u00f3ricas, aulas pr\u00e1ticas, aulas de exerc\u00edcios. MANTIDO\" `\n    \"Provas, atividades em grupo e atividades individuais.\"\n\nReplace-InRange $avaliacao `\n    \"M = (0,8P + 0,2T)P = m\u00e9dia aritm\u00e9tica de duas provas escritasT = M\u00e9dia das notas de trabalhos e exerc\u00edciosM = M\u00e9dia de aproveitamento do alunoAprova\u00e7\u00e3o com m\u00e9dia de aproveitamento maior ou igual a 5,0 e no m\u00ednimo 70% de frequ\u00eancia \u00e0s aulas.\" `\n    \"M\u00e9dia das atividades avaliativas\"\n\nReplace-InRange $avaliacao `\n    \"MF = (0,5 M + 0,5 R)M = M\u00e9dia de aproveitamento do aluno, antes da recupera\u00e7\u00e3oR = Nota de uma prova de recupera\u00e7\u00e3oMF = nota final de aproveitamento, ap\u00f3s a recupera\u00e7\u00e3oAprova\u00e7\u00e3o com m\u00e9dia final de aproveitamento maior ou igual a 5,0.A recupera\u00e7\u00e3o dever\u00e1 consistir de uma prova escrita englobando a mat\u00e9ria toda do semestre.Ter\u00e1 direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\" `\n    \"MF = (0,5 M + 0,5 R) M = M\u00e9dia de aproveitamento do aluno, antes da recupera\u00e7\u00e3o R = Nota de uma prova de recupera\u00e7\u00e3o MF = nota final de aproveitamento, ap\u00f3s a recupera\u00e7\u00e3o Aprova\u00e7\u00e3o com m\u00e9dia final de aproveitamento maior ou igual a 5,0. A recupera\u00e7\u00e3o dever\u00e1 consistir em uma prova escrita englobando a mat\u00e9ria toda do semestre. Ter\u00e1 direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\"\n"}
